$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the title/metadata string in A1
$ws.Range("A1").Value = "#Title=TestArrayTransposed - HeaderDepth=3 - IsTransposed=True - DateTime=2025-09-18 09:28:30.340200 - DatamodelUrl=None"

# Update date/datetime/time serial values in row 10-12, columns D:F
$ws.Range("D10:F10").Value = 45918
$ws.Range("D11:F11").Value = 45918.39477812201
$ws.Range("D12:F12").Value = 0.3947781221296296
